$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.036.18"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.69"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.657.80"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.362.05"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.930.14"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.98"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.77"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.015.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.08"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.526.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -1.23%  "
